$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" values in this sheet are stored as text (inline strings),
# not numbers, even though they look numeric. Force the NumberFormat to
# Text ("@") before assigning so Excel keeps the assigned value as text
# instead of silently re-casting it to a numeric value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "248.66"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.64"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.262"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.338"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8053"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8993"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1402"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07448"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03089"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03005"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09384"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.859"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001583"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04766"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.01828"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0005801"
$ws.Range("E19").Value = "18OneONE"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006451"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004990"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.001000"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0001500"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.695"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.197"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1307"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03968"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006830"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1067"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002768"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007715"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005594"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4990"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2080"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01010"
